{"js": "// Highlight the \"Queues:\" and \"Applications of Stacks and Queues:\" bullet\n// paragraphs (under \"8. Stacks and Queues\") with a cyan highlight, matching\n// the existing cyan highlight already applied to the \"Stacks:\" bullet above\n// them. This sets the highlight on the paragraph mark as well as on every\n// run in the paragraph (mirroring the sibling \"Stacks:\" paragraph).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// The two target bullet paragraphs are identified by the text their runs\n// start with (the bold \"label\" run followed by the description run).\nconst targets = [\"Queues:\", \"Applications of Stacks and Queues:\"];\n\nfor (const para of paragraphs.items) {\n  const text = para.text || \"\";\n  if (targets.some((t) => text.indexOf(t) === 0)) {\n    // Setting highlightColor on the paragraph's Font applies it to the\n    // paragraph mark (w:pPr/w:rPr) and to every run currently in the\n    // paragraph (w:r/w:rPr) \u2014 exactly mirroring the sibling \"Stacks:\" bullet.\n    para.font.highlightColor = \"cyan\";\n  }\n}\n\nawait context.sync();\n", "ps1": "# Highlight the \"Queues:\" and \"Applications of Stacks and Queues:\" bullet\n# paragraphs (under \"8. Stacks and Queues\") with a cyan highlight, matching\n# the existing cyan highlight already applied to the \"Stacks:\" bullet above\n# them. Setting Font.HighlightColorIndex on the paragraph's Range applies the\n# highlight to the paragraph mark (w:pPr/w:rPr) as well as to every run in\n# the paragraph (w:r/w:rPr) \u2014 mirroring the sibling \"Stacks:\" paragraph.\n\n$d = $word.ActiveDocument\n\n$targets = @(\"Queues:\", \"Applications of Stacks and Queues:\")\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    foreach ($t in $targets) {\n        if ($text.StartsWith($t)) {\n            $p.Range.Font.HighlightColorIndex = \"wdTurquoise\"\n            break\n        }\n    }\n}\n"}
